$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Fill Rate" values for each flower section (column C).
# The dependent "Watering Rate" formulas in column C (B/C) recalc automatically.
$ws.Range("C5").Value = 50
$ws.Range("C8").Value = 45
$ws.Range("C11").Value = 40
$ws.Range("C14").Value = 35
$ws.Range("C17").Value = 35
$ws.Range("C20").Value = 30

$excel.CalculateFullRebuild()

# Move the selected/active cell to D22, matching the author's final selection.
$ws.Range("D22").Select()
